# Applies a 3-row rotation of observation records (rows 2-4)
# Row 2 <= old Row 3 data, Row 3 <= old Row 4 data (with biotope-description update
# and substrate fields cleared), Row 4 <= old Row 2 data (with substrate fields set).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 111697636
$ws.Range("B2").Value = 88489
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 1962
$ws.Range("F2").Value = 'Vaddporing'
$ws.Range("G2").Value = 'Anomoporia kamtschatica'
$ws.Range("H2").Value = '(Parmasto) Bondartseva'
$ws.Range("J2").Value = 'fruktkroppar'
$ws.Range("M2").Value = ""
$ws.Range("Q2").Value = 373112.5181173298
$ws.Range("R2").Value = 6865358.590016441
$ws.Range("Z2").Value = '19:00'
$ws.Range("AB2").Value = '19:00'
$ws.Range("AC2").Value = 'Växer under rötad gammal silverved'
$ws.Range("AI2").Value = 'Kontinuitetsskog. Tallskog'

# --- Row 3 ---
$ws.Range("A3").Value = 111697304
$ws.Range("B3").Value = 8377
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 106545
$ws.Range("F3").Value = 'Mindre märgborre'
$ws.Range("G3").Value = 'Tomicus minor'
$ws.Range("H3").Value = '(Hartig, 1834)'
$ws.Range("J3").Value = ""
$ws.Range("M3").Value = 'färska gnagspår'
$ws.Range("Q3").Value = 373090.8741807578
$ws.Range("R3").Value = 6865424.499624529
$ws.Range("AC3").Value = ""
$ws.Range("AI3").Value = 'Luckig tallskog. K-skog'
$ws.Range("AJ3").Value = ""
$ws.Range("AK3").Value = ""
$ws.Range("AO3").Value = ""

# --- Row 4 ---
$ws.Range("A4").Value = 111697236
$ws.Range("Q4").Value = 373121.3523494597
$ws.Range("R4").Value = 6865443.651501717
$ws.Range("Z4").Value = '00:00'
$ws.Range("AB4").Value = '00:00'
$ws.Range("AI4").Value = 'Tallskog. Kontinuitetsskog'
$ws.Range("AJ4").Value = 'tall'
$ws.Range("AK4").Value = 'Pinus sylvestris'
$ws.Range("AO4").Value = 'Pinus sylvestris'
